# Applies the "add content to index and java_lib" edit to Sheet1.
# New rows are appended for: jdk (A13), j2ee (row 14, Servlet example),
# and two GSON/Gson helper rows (rows 15-16). Row heights for rows 9/10
# are widened (their C-column code blocks now need more vertical room),
# and the active selection moves to C15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 13: add the "jdk" language label; the existing Regex/Pattern
#     code sample in B13/C13 is left untouched. ---
$ws.Range("A13").Value = "jdk"

# --- Row 14 (new): j2ee / HttpServletResponse Header Setting / servlet sample ---
$ws.Range("A14").Value = "j2ee"
$ws.Range("B14").Value = "HttpServletResponse Header Setting"

$servletCode = @'
import java.io.*;
import javax.servlet.*;
import javax.servlet.http.*;
public class ExampServlet extends HttpServlet {
  public void doPost(HttpServletRequest request, 
         HttpServletResponse response)
        throws ServletException, IOException
  {
    response.setContentType("text/html");
    PrintWriter out = response.getWriter();
    out.println("<title>Example</title>" +
       "<body bgcolor=FFFFFF>");
    out.println("<h2>Button Clicked</h2>");
    String DATA = request.getParameter("DATA");
    if(DATA != null){
      out.println(DATA);
    } else {
      out.println("No text entered.");
    }
    out.println("<P>Return to 
        <A HREF="../simpleHTML.html">Form</A>");
    out.close();
  }
}
'@
$ws.Range("C14").Value = $servletCode

# --- Row 15 (new): GSON / map json to class object / fromJson ---
$ws.Range("A15").Value = "GSON"
$ws.Range("B15").Value = "map json to class object"
$ws.Range("C15").Value = "new Gson().fromJson(jsonString, ClassOfObject.class);"

# --- Row 16 (new): GSON / parse class object to json(string) / toJson ---
$ws.Range("A16").Value = "GSON"
$ws.Range("B16").Value = "parse olass object to json(string)"
$ws.Range("C16").Value = "new Gson().toJson(obj);"

# --- Row heights: rows 9 & 10 now hold taller code blocks in column C ---
$ws.Rows.Item(9).RowHeight = 316.5
$ws.Rows.Item(10).RowHeight = 85.5

# keep the new rows at the standard 33pt row height
$ws.Rows.Item(13).RowHeight = 33
$ws.Rows.Item(14).RowHeight = 33
$ws.Rows.Item(15).RowHeight = 33
$ws.Rows.Item(16).RowHeight = 33

# --- Selection moves to the last-edited cell ---
$ws.Range("C15").Select() | Out-Null
